# Model 33 / M756 Tanks workbook fix-ups:
#  1) Tanks: add missing "M" prefix to M756 top-level item number
#     -> handled generically below since the top level item string is
#        shared across both sheets (Materials!A, Operations!A, Materials!I10).
#  2) Model 33: top level item is now defined as "M33M" instead of the old
#     drawing number "8027958".
#  3) Model 33: remove the self-referencing top-level-assembly component row
#     (old row 10 in Materials), shifting the remaining component rows up by
#     one and dropping the now-duplicate last row.

$wb = $excel.ActiveWorkbook
$wsMat = $wb.Worksheets.Item("Materials")
$wsOps = $wb.Worksheets.Item("Operations")

# --- Step 1: remove the self-reference -----------------------------------
# Rows 10-67 of the Materials sheet are the BOM lines for the Model 33
# top-level assembly. Row 10 used to reference the assembly's own drawing
# number/description (a self-reference). Drop it by shifting every other
# line's Material / Material Description / Quantity / U-M up one row, then
# deleting the now-redundant last row (67).
#
# Use Copy + PasteSpecial(values) rather than a plain Value2 assignment for
# the Material column (I) because several Material numbers are purely
# numeric-looking text (e.g. "8034376"); a raw Value2 round-trip would let
# Excel re-interpret them as numbers and silently change their stored type.
$wsMat.Range("I11:J67").Copy()
$wsMat.Range("I10:J66").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Quantity (S) and U/M (U) are never ambiguous (S is genuinely numeric, U is
# short text like "EA"/"IN"), so a plain array copy is fine and keeps things
# simple. T (the "Unit" label column) rides along unchanged since it is
# constant for every row in this block.
$wsMat.Range("S10:U66").Value2 = $wsMat.Range("S11:U67").Value2()

# The old last row (67) is now a duplicate of row 66 and gets removed
# entirely, shifting the sheet's used range from A1:AQ67 to A1:AQ66.
$wsMat.Rows("67").Delete()

# --- Step 2: rename the top level item ------------------------------------
# The top-level Model 33 assembly item is now identified as "M33M" instead
# of its old drawing number "8027958". That value is shown in the Item
# column (A) of every Materials BOM line for this assembly, and in the
# Operations sheet's Item column (A) for its routing.
$wsMat.Range("A10:A66").Value2 = "M33M"
$wsOps.Range("A25").Value2 = "M33M"
